$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.02286966666666667
$ws.Range("H2").Value = 0.068609
$ws.Range("M2").Value = 1.041066666666667
$ws.Range("N2").Value = 3.1232
$ws.Range("O2").Value = 0.06773090411171609
$ws.Range("P2").Value = 0.06773090411171608
$ws.Range("Q2").Value = 0.02380884764444444
$ws.Range("R2").Value = 0.2142796288
$ws.Range("S2").Value = 0.06773090411171609
$ws.Range("T2").Value = 0.06773090411171608

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.02286966666666667
$ws.Range("H3").Value = 0.068609
$ws.Range("M3").Value = 7.745649666666666
$ws.Range("O3").Value = 0.5039253216469766
$ws.Range("P3").Value = 0.5039253216469765
$ws.Range("Q3").Value = 0.1771404259934444
$ws.Range("R3").Value = 1.594263833941
$ws.Range("S3").Value = 0.5039253216469766
$ws.Range("T3").Value = 0.5039253216469765

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.02286966666666667
$ws.Range("H4").Value = 0.068609
$ws.Range("M4").Value = 6.583913666666667
$ws.Range("O4").Value = 0.4283437742413074
$ws.Range("P4").Value = 0.4283437742413074
$ws.Range("Q4").Value = 0.1505719109187778
$ws.Range("R4").Value = 1.355147198269
$ws.Range("S4").Value = 0.4283437742413074
$ws.Range("T4").Value = 0.4283437742413074
